# Ticket 35: POSMasterfile template - add "POS Name" column.
#
# A new column is inserted at C (pushing the former Category/SubCategory/
# SRP/DeliveryPrice/TableVibePrice/Active columns one to the right), headed
# "POS Name", and populated per-row with the same text already present in
# column B (Item Description) - i.e. the newly-added column is "completed"
# by copying the description down as the POS name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column C (Category). This shifts
# Category/SubCategory/SRP/DeliveryPrice/TableVibePrice/Active from C:G to
# D:H (and Active from H to I), carrying their styles/widths along.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "POS Name"

# New column width matches column B's (the header text's own column),
# as a fixed custom width (no "best fit" autofit flag).
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Fill the new column's data rows with the Item Description value that is
# already in column B for that row.
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Value2
}

# Leave the just-completed column selected, matching the author's
# post-edit selection.
$ws.Range("C2:C7").Select()
